$wb = $excel.ActiveWorkbook

# Sheet4: currently the active/selected sheet with selection A2.
# After the edit it is no longer the active sheet, and its stored
# selection becomes A25:B25.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A25:B25").Select()

# Sheet5: gains a new row of data (row 15) and becomes the active sheet
# with selection A15:B15.
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("A15").Value = "李四12313131"
$ws5.Range("B15").Value = -200

$ws5.Activate()
$ws5.Range("A15:B15").Select()
